$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = @(0.4141732146861696, 1.605476919376253, 4.285373823189503, 2.070114446881984, 2.050671332274695, 46)
    3  = @(-0.01716789264613588, 1.651420560915961, 4.541767649565558, 2.131142334421978, 2.155153870726067, 45)
    4  = @(0.4406255020385922, 1.625614656034232, 4.490980664127665, 2.119193399415841, 2.096844320931413, 44)
    5  = @(0.03988787390244446, 1.807057129213736, 5.145130224049834, 2.268287949985591, 2.294777636324286, 43)
    6  = @(0.455682680408077, 1.62477962444376, 4.358939589200301, 2.087807364006627, 2.062169714142621, 42)
    7  = @(0.1777087621310496, 1.698096300376475, 4.364176648270135, 2.08906118825422, 2.107346945862329, 41)
    8  = @(0.5039136335280825, 1.664940487861543, 4.498869753835015, 2.121053925253909, 2.086572653845029, 40)
    9  = @(0.1430518963212211, 1.71910815220499, 4.463347462520725, 2.112663594262164, 2.135369196980537, 39)
    10 = @(0.5887292808226858, 1.678991187388869, 4.640656357710296, 2.154218270675071, 2.100026121338619, 38)
    11 = @(0.1084470123171574, 1.756824810856007, 4.640724073720614, 2.15423398769043, 2.181179871918736, 37)
}

foreach ($row in $newValues.Keys) {
    $vals = $newValues[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 6).Value = $vals[4]
    $ws.Cells.Item($row, 7).Value = $vals[5]
}
